$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the two new columns, matching the existing header style
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for I2:I60 and J2:J60
$iValues = @(7,7,12,8,8,8,7,9,9,10,7,6,6,9,8,7,9,8,6,7,9,8,7,7,7,7,4,4,6,8,7,5,9,8,6,10,5,7,8,7,8,8,4,5,7,8,7,7,6,6,6,7,8,6,8,6,5,6,5)
$jValues = @(8,7,12,8,8,8,7,9,9,10,7,6,6,9,8,7,9,8,7,8,9,9,7,8,8,7,5,5,6,8,7,6,9,8,7,10,6,8,8,7,8,8,5,6,8,8,7,7,6,6,7,7,8,7,8,7,5,6,5)

for ($r = 2; $r -le 60; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
